# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet at the front of the workbook with the
#    player's basic info (id, name, batting hand, bowling style).
# 2. Rename the MATCH_CARD_LINK column (full howstat.com URL) to MATCH_CODE
#    on both the "ODI Batting" and "ODI Bowling" sheets, and replace the
#    stored values with just the numeric match code that used to be the
#    query-string parameter of that URL.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text even when the value looks like a
    # number (e.g. "5937", "4465") - matches the original workbook, which
    # stores every cell (headers & data alike) as inline/shared strings
    # rather than numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop back to the default "Normal" style so no stray number-format
    # style is left attached to the cell (mirrors the unstyled data cells
    # already present in the source file).
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the current first sheet.
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Style the header row like the existing header rows on the other sheets
# (bold, bordered, centered / top aligned).
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1        # xlContinuous

Set-TextValue $playerInfo.Range("A2") "5937"
$playerInfo.Range("B2").Value = "Wanigamuni Ramesh Tarinda Mendis"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ---------------------------------------------------------------------
# 2. ODI Batting sheet - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

Set-TextValue $batting.Range("D2") "4465"
Set-TextValue $batting.Range("D3") "4469"
Set-TextValue $batting.Range("D4") "4485"
Set-TextValue $batting.Range("D5") "4527"

# ---------------------------------------------------------------------
# 3. ODI Bowling sheet - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

Set-TextValue $bowling.Range("B2") "4465"
Set-TextValue $bowling.Range("B3") "4485"
Set-TextValue $bowling.Range("B4") "4527"
